$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the styled column-A / header formatting (style index "1": bold, bordered,
# centered) down to the two new rows (16 and 17) by copying the format from the
# last existing styled row before rewriting the table contents.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Clear all the existing values in the table (header + data rows) so the whole
# table (including the two brand new rows) can be rewritten in one consistent
# pass, top to bottom, left to right, matching the final layout exactly.
$ws.Range("A1:E17").ClearContents()

# Header row
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "from_bus"
$ws.Range("D1").Value = "to_bus"
$ws.Range("E1").Value = "in_service"

# Data rows 2-17
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "line1"
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = $true

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "line2"
$ws.Range("C3").Value = 9
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "line3"
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = $true

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "line4"
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = $true

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "line5"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = $false

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "line6"
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = $true

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
